# This workbook's data rows (2-41) get reshuffled: the values in columns
# D, L, M, N, O, P, R and S move to different rows while columns
# A, B, C, E, F, G, H, I, J, K, Q, T stay constant for every row and are
# left untouched. The mapping below gives, for each target row (2..41 in
# order), which row currently holds the data that should end up there.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 41

# sourceForTarget[i] = source row whose D/L/M/N/O/P/R/S values should be
# written into target row (firstRow + i)
$sourceForTarget = @(34,18,24,23,37,20,29,13,25,17,5,15,2,28,6,35,19,27,10,30,40,7,12,16,36,21,26,9,11,31,8,39,33,14,3,38,41,32,22,4)

$cols = @("D","L","M","N","O","P","R","S")

# Snapshot the current values for the columns that move, before we start
# overwriting any of them.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Write the shuffled values back out row by row.
for ($i = 0; $i -lt $sourceForTarget.Length; $i++) {
    $targetRow = $firstRow + $i
    $srcRow = $sourceForTarget[$i]
    $srcVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$targetRow").Value2 = $srcVals[$col]
    }
}
